# Apply the "Mise à jour du fichier Excel public" update:
#  1. "Prix Spot" sheet: column N (27-jun), rows 2-25 had numeric values that are
#     now unavailable -> replace with the text placeholder "-".
#  2. "CO2" sheet: insert two new daily rows (2025-06-21 and 2025-06-22) before the
#     existing 2025-06-23 row, pushing the final 2025-06-25 row down.

$wb = $excel.ActiveWorkbook

# --- 1. "Prix Spot": N2:N25 -> "-" -------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Range("N2:N25").Value = "-"

# --- 2. "CO2": insert two rows with new dates ---------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

# Push the current row 7 (2025-06-23 / 71.88) and row 8 (2025-06-25 / 70.17) down
# by two rows so we can insert 2025-06-21 and 2025-06-22 ahead of them.
$wsCo2.Rows.Item(8).Insert()
$wsCo2.Rows.Item(8).Insert()

# Force the date column to be written back as plain text (matches the existing
# text-typed date cells in the sheet) instead of being auto-parsed into a date
# serial value, then restore the default "Normal" style so no stray number
# format lingers on the cells.
$wsCo2.Range("A7:A9").NumberFormat = "@"

$wsCo2.Range("A7").Value = "2025-06-21"
$wsCo2.Range("B7").Value = 72.2

$wsCo2.Range("A8").Value = "2025-06-22"
$wsCo2.Range("B8").Value = 72.2

$wsCo2.Range("A9").Value = "2025-06-23"
$wsCo2.Range("B9").Value = 71.88

$wsCo2.Range("A7:A9").Style = "Normal"
